# Apply the "zero_before_threshold" recalculation results to the
# Step3_DataPts_* sheets: First_Noticeable_Increase_Index (C),
# First_Noticeable_Increase_Cumulative_Value (E), and Pulse_Width (G)
# shift for every signal segment row (2-6) on every threshold sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws.Range("C2").Value = 87
$ws.Range("E2").Value = 0.04550896254114215
$ws.Range("G2").Value = 15
$ws.Range("C3").Value = 87
$ws.Range("E3").Value = 0.03016927218527709
$ws.Range("G3").Value = 15
$ws.Range("C4").Value = 88
$ws.Range("E4").Value = 0.02970608908351221
$ws.Range("G4").Value = 8
$ws.Range("C5").Value = 88
$ws.Range("E5").Value = 0.01040862661396182
$ws.Range("G5").Value = 8
$ws.Range("C6").Value = 87
$ws.Range("E6").Value = 0.03924871506819344
$ws.Range("G6").Value = 9

$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws.Range("C2").Value = 87
$ws.Range("E2").Value = 0.04550896254114215
$ws.Range("G2").Value = 35
$ws.Range("C3").Value = 87
$ws.Range("E3").Value = 0.03016927218527709
$ws.Range("G3").Value = 36
$ws.Range("C4").Value = 88
$ws.Range("E4").Value = 0.02970608908351221
$ws.Range("G4").Value = 34
$ws.Range("C5").Value = 88
$ws.Range("E5").Value = 0.01040862661396182
$ws.Range("G5").Value = 37
$ws.Range("C6").Value = 87
$ws.Range("E6").Value = 0.03924871506819344
$ws.Range("G6").Value = 35

$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws.Range("C2").Value = 87
$ws.Range("E2").Value = 0.04550896254114215
$ws.Range("G2").Value = 70
$ws.Range("C3").Value = 87
$ws.Range("E3").Value = 0.03016927218527709
$ws.Range("G3").Value = 48
$ws.Range("C4").Value = 88
$ws.Range("E4").Value = 0.02970608908351221
$ws.Range("G4").Value = 47
$ws.Range("C5").Value = 88
$ws.Range("E5").Value = 0.01040862661396182
$ws.Range("G5").Value = 50
$ws.Range("C6").Value = 87
$ws.Range("E6").Value = 0.03924871506819344
$ws.Range("G6").Value = 48

$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws.Range("C2").Value = 87
$ws.Range("E2").Value = 0.04550896254114215
$ws.Range("G2").Value = 92
$ws.Range("C3").Value = 87
$ws.Range("E3").Value = 0.03016927218527709
$ws.Range("G3").Value = 76
$ws.Range("C4").Value = 88
$ws.Range("E4").Value = 0.02970608908351221
$ws.Range("G4").Value = 79
$ws.Range("C5").Value = 88
$ws.Range("E5").Value = 0.01040862661396182
$ws.Range("G5").Value = 78
$ws.Range("C6").Value = 87
$ws.Range("E6").Value = 0.03924871506819344
$ws.Range("G6").Value = 79
